$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.724.50'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '1.538.50'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '289.48'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('D7').Value = '0.3892'
$ws.Range('E7').Value = '  +3.22%  '
$ws.Range('D8').Value = '0.3185'
$ws.Range('E8').Value = '  -1.77%  '
$ws.Range('D9').Value = '42.93'
$ws.Range('E9').Value = '  +3.60%  '
$ws.Range('D10').Value = '0.07201'
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('D11').Value = '1.057'
$ws.Range('E11').Value = '  -6.43%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '5.635'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').Value = '18.58'
$ws.Range('E14').Value = '  -5.56%  '
$ws.Range('D15').Value = '6.610'
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.543.41'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001111'
$ws.Range('E17').Value = '  +2.58%  '
$ws.Range('D18').Value = '0.06585'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = '83.17'
$ws.Range('E19').Value = '  -2.28%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '6.141'
$ws.Range('E21').Value = '  -4.93%  '
$ws.Range('D22').Value = '15.38'
$ws.Range('E22').Value = '  -3.81%  '
$ws.Range('D23').Value = '10.85'
$ws.Range('E23').Value = '  -6.33%  '
$ws.Range('D24').Value = '2.407'
$ws.Range('E24').Value = '  +7.21%  '
$ws.Range('D25').Value = '21.727.74'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = '2.369'
$ws.Range('E26').Value = '  -6.55%  '
$ws.Range('D27').Value = '146.14'
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('D28').Value = '18.36'
$ws.Range('E28').Value = '  -3.15%  '
$ws.Range('D29').Value = '4.834'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '1.715.80'
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').Value = '117.61'
$ws.Range('E31').Value = '  -2.22%  '
$ws.Range('D32').Value = '0.9698'
$ws.Range('E32').Value = '  -13.05%  '
$ws.Range('D33').Value = '5.905'
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('D34').Value = '0.08192'
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('D35').Value = '8.809'
$ws.Range('E35').Value = '  -5.16%  '
$ws.Range('D36').Value = '0.06087'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('D37').Value = '5.134'
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('D38').Value = '1.488'
$ws.Range('E38').Value = '  -7.12%  '
$ws.Range('D39').Value = '0.02199'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('D40').Value = '0.2037'
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('D41').Value = '1.190'
$ws.Range('E41').Value = '  -2.12%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '10.66'
$ws.Range('E43').Value = '  -2.41%  '
$ws.Range('D44').Value = '0.5739'
$ws.Range('E44').Value = '  -3.55%  '
$ws.Range('D45').Value = '13.08'
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('D46').Value = '3.740'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('D47').Value = '0.5502'
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').Value = '1.167'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '117.41'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '1.866'
$ws.Range('E50').Value = '  -4.04%  '
$ws.Range('D51').Value = '0.06724'
$ws.Range('E51').Value = '  -3.11%  '
